$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "[Object Detection] YOLO Define Optimal Anchor Box :: YOLO v5, YOLO v6 autoanchor"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Object-Detection-YOLO-Optimal-Anchor-Box-YOLO-v5-YOLO-v6-autoanchor"

$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"

$ws.Range("D32").Value = "불균형 데이터분석을 위한 샘플링 기법"
$ws.Range("E32").Value = "https://dodonam.tistory.com/385"

$ws.Range("D46").Value = "Subdural hemorrhage (SDH; 경막하출혈)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/483"

$ws.Range("D51").Value = "[flask] Ajax로 json 데이터 post 전달시 flask 서버 get_json() 관련 오류 해결 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/flask-Ajax%EB%A1%9C-json-%EB%8D%B0%EC%9D%B4%ED%84%B0-post-%EC%A0%84%EB%8B%AC%EC%8B%9C-flask-%EC%84%9C%EB%B2%84-getjson-%EA%B4%80%EB%A0%A8-%EC%98%A4%EB%A5%98-%ED%95%B4%EA%B2%B0-%EB%B0%A9%EB%B2%95"
